$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
# D1 label renamed from years_offered -> acad_years_offered
$ws.Range('D1').Value = 'acad_years_offered'

# New column F ("Course Description") with header styled like the other headers
$ws.Range('E1').Copy()
$ws.Range('F1').PasteSpecial(-4122)
$ws.Range('F1').Value = 'Course Description'

# --- Per-row data (row, acad_years_offered, OFS_CBI, Course Description) ---
# Kept as parallel flat arrays (nested @() literals get flattened by this host).
$rowNums = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 82, 83, 84, 85, 86, 87, 88, 89)
$dVals = @(2, 2, 2, 2, 4, 4, 2, 1, 3, 4, 4, 3, 2, 4, 4, 4, 4, 4, 4, 2, 4, 2, 4, 4, 4, 4, 3, 4, 3, 2, 4, 4, 4, 1, 4, 4, 3, 4, 3, 3, 2, 4, 4, 4, 3, 2, 4, 4, 4, 3, 2, 4, 4, 4, 4, 4, 4, 4, 3, 1, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 1, 4, 4, 4, 4, 4, 4, 3, 3, 3, 4, 3, 4, 4, 4, 2, 4, 2)
$eVals = @(3, 3, 3, 3, 2, 2, 3, 3, 3, 1, 1, 3, 3, 2, 1, 1, 1, 2, 1, 3, 1, 3, 2, 2, 1, 2, 3, 2, 3, 3, 2, 1, 1, 3, 2, 2, 3, 2, 3, 3, 3, 2, 1, 2, 3, 3, 1, 1, 1, 3, 3, 2, 2, 2, 2, 2, 2, 2, 3, 3, 2, 2, 1, 2, 2, 2, 2, 1, 1, 2, 3, 2, 2, 2, 2, 2, 2, 3, 3, 3, 2, 3, 1, 2, 1, 3, 2, 3)
$fVals = @('Applied Cyberinfrastruct Conc', 'Game AI', 'Designing an Installation', 'Data Ethics', 'Foundations of Information', 'Information Research Methods', 'Bayesian Modeling & Inference', 'Foundations of Data Science', 'Computational Social Science', 'Organization/Information', 'Intro: Human Computer Interact', 'Intro to Digital Cultures', 'Info Trust and Manipulation', 'Ethical Issues in Information', 'Intro to Machine Learning', 'Data Mining/Discovery', 'Virtual Reality', 'Algorithms for Games', 'Data Analysis and Visualizatio', 'Applied Cyberinfrastruct Conc', 'Data Warehousing in the Cloud', 'Data Science, Public Interests', 'Introduction To Archives', 'Artificial Intelligence', 'Game Development', 'Advanced Game Development', 'Applied NLP', 'Neural Networks', 'Social Justice in Info Service', 'STEM Games', 'Leadership & the Info Org', 'Database Dev And Mgmt', 'Intro Info Technology', 'User Interf+Website Dsgn', 'Information Security', 'Science Information', 'SQL/NoSQL Databases', 'Data for the Semantic Web', 'Intellectual Property/Copyrigh', 'Managing the Information Org', 'Adv ML Apps', 'Adv Archives: Apprsl & Dscr', 'Intro Digital Curation/Preserv', 'Intro Applied Technology', 'Managing Digital Info', 'Adv Digital Collections', 'Found Libr+Info Services', 'Rsrch Mth/Libr+Info Prof', 'Organization/Information', 'Intro to Digital Cultures', 'Info Trust and Manipulation', 'Ethical Issues in Information', 'Chldrn+Young Adult Lit', 'Early Chlhd+Public Librs', 'Cataloging+Metadata Mgmt', 'Info Intermediation', 'Business Information', 'Introduction To Archives', 'Preservation', 'Curating & Preserving Media', 'Inf Env/Non-dominant Pers', 'Documnt Divrs Cult+Comms', 'Social Justice in Info Service', 'Mktng Library+Info Srvcs', 'Collection Management', 'Rdrs Advisory/Publ Libr', 'Leadership & the Info Org', 'Database Dev And Mgmt', 'Intro Info Technology', 'Government Information', 'User Interf+Website Dsgn', 'Science Information', 'Data for the Semantic Web', 'How to Teach Info Literacy', 'Young Adults+Public Libr', 'eLearning for Librarians', 'Intellectual Property/Copyrigh', 'Learning Design Lib/Info Prof', 'Special Topics in LIS', 'Managing the Information Org', 'Adv Archives: Apprsl & Dscr', 'Comm focused Archives/Museums', 'Intro Digital Curation/Preserv', 'Intro Applied Technology', 'Managing Digital Info', 'Adv Digital Collections', 'Digital Info Mgmt Capstn', 'Applied Cyberinfrastruct Conc')

for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = [int]$rowNums[$i]
    $d = [int]$dVals[$i]
    $e = [int]$eVals[$i]
    $f = [string]$fVals[$i]
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
}
